$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# Row 3 (C3): "Hypotheses 1 & 2" -> "Hypotheses 1 " (drop "& 2")
$ws.Range("C3").Characters(22, 3).Text = ""

# Row 5 (C5): "...hypothesis tests 1 & 2" -> "...hypothesis tests 1 " (drop "& 2")
$ws.Range("C5").Characters(57, 3).Text = ""

# Row 13 (C13): "Hypothesis 5" -> "Hypothesis 2"
$ws.Range("C13").Characters(20, 1).Text = "2"

# Row 15 (C15): "Data Analysis for hypothesis test 5" -> "Data Analysis for hypothesis test 2"
$ws.Range("C15").Characters(43, 1).Text = "2"

# Update the active selection to match the new state
$ws.Range("C15").Select()
